$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1003.375
$ws.Range("J17").Value = 1003.375
$ws.Range("L17").Value = 3010.125
$ws.Range("N17").Value = -3346.125
$ws.Range("H64").Value = 4699.5
$ws.Range("I64").Value = 3865.6667
$ws.Range("J64").Value = 5199.8
$ws.Range("K64").Value = 3865.6667
$ws.Range("L64").Value = 5199.8
$ws.Range("M64").Value = -3617.6667
$ws.Range("N64").Value = -5695.8
$ws.Range("H67").Value = 4699.5
$ws.Range("I67").Value = 3865.6667
$ws.Range("J67").Value = 5199.8
$ws.Range("K67").Value = 3865.6667
$ws.Range("L67").Value = 5199.8
$ws.Range("M67").Value = -3007.6667
$ws.Range("N67").Value = -6915.8
$ws.Range("H98").Value = 576.7273
$ws.Range("I98").Value = 555.5
$ws.Range("J98").Value = 633.3333
$ws.Range("K98").Value = 555.5
$ws.Range("L98").Value = 633.3333
$ws.Range("M98").Value = 942.5
$ws.Range("N98").Value = -3629.3333
$ws.Range("H99").Value = 573
$ws.Range("J99").Value = 1999
$ws.Range("L99").Value = 5997
$ws.Range("N99").Value = -8993
$ws.Range("H101").Value = 33339300
$ws.Range("I101").Value = 50008450
$ws.Range("J101").Value = 999
$ws.Range("K101").Value = 150025350
$ws.Range("L101").Value = 2997
$ws.Range("M101").Value = -150023728
$ws.Range("N101").Value = -6241
$ws.Range("H113").Value = 4498.5
$ws.Range("I113").Value = 4498.5
$ws.Range("K113").Value = 4498.5
$ws.Range("M113").Value = -1244.5
$ws.Range("H116").Value = 14337
$ws.Range("J116").Value = 6503
$ws.Range("L116").Value = 6503
$ws.Range("N116").Value = -13387
$ws.Range("H122").Value = 576.7273
$ws.Range("I122").Value = 555.5
$ws.Range("J122").Value = 633.3333
$ws.Range("K122").Value = 1666.5
$ws.Range("L122").Value = 1899.9999
$ws.Range("M122").Value = 783.5
$ws.Range("N122").Value = -6799.9999
$ws.Range("H125").Value = 12549.833
$ws.Range("J125").Value = 18100
$ws.Range("L125").Value = 162900
$ws.Range("N125").Value = -167820
$ws.Range("H138").Value = 3943.8333
$ws.Range("I138").Value = 3999.8333
$ws.Range("J138").Value = 3915.8333
$ws.Range("K138").Value = 11999.4999
$ws.Range("L138").Value = 11747.4999
$ws.Range("M138").Value = -6859.499899999999
$ws.Range("N138").Value = -22027.4999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 150000
$ws.Range("J64").Value = 150000
$ws.Range("L64").Value = 150000
$ws.Range("N64").Value = -150496
$ws.Range("H67").Value = 150000
$ws.Range("J67").Value = 150000
$ws.Range("L67").Value = 150000
$ws.Range("N67").Value = -151716
$ws.Range("H97").Value = 727.2222
$ws.Range("I97").Value = 762.1429000000001
$ws.Range("K97").Value = 762.1429000000001
$ws.Range("M97").Value = -266.1429000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2785.5715
$ws.Range("I105").Value = 2499.8
$ws.Range("K105").Value = 2499.8
$ws.Range("M105").Value = -752.8000000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2761.8572
$ws.Range("I2").Value = 781.3333
$ws.Range("K2").Value = 781.3333
$ws.Range("M2").Value = -668.3333
$ws.Range("H31").Value = 2963.8333
$ws.Range("I31").Value = 2756.6
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 2756.6
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -2461.6
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 2963.8333
$ws.Range("I34").Value = 2756.6
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 2756.6
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -2554.6
$ws.Range("N34").Value = -4404
$ws.Range("H86").Value = 10340.5
$ws.Range("I86").Value = 10962.571
$ws.Range("K86").Value = 10962.571
$ws.Range("M86").Value = -9839.571
$ws.Range("H89").Value = 10340.5
$ws.Range("I89").Value = 10962.571
$ws.Range("K89").Value = 54812.855
$ws.Range("M89").Value = -49196.855
$ws.Range("H122").Value = 3050
$ws.Range("I122").Value = 1466.6666
$ws.Range("J122").Value = 4633.3335
$ws.Range("K122").Value = 4399.9998
$ws.Range("L122").Value = 13900.0005
$ws.Range("M122").Value = -1949.9998
$ws.Range("N122").Value = -18800.0005
$ws.Range("H132").Value = 9333
$ws.Range("I132").Value = 9333
$ws.Range("K132").Value = 27999
$ws.Range("M132").Value = -25469
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.05263
$ws.Range("I2").Value = 20.363636
$ws.Range("J2").Value = 52.875
$ws.Range("K2").Value = 122.181816
$ws.Range("L2").Value = 317.25
$ws.Range("M2").Value = -9.181815999999998
$ws.Range("N2").Value = -543.25
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6338
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6204
$ws.Range("H119").Value = 3499
$ws.Range("I119").Value = 3499
$ws.Range("K119").Value = 10497
$ws.Range("M119").Value = -5659
$ws.Range("H122").Value = 150
$ws.Range("I122").Value = 150
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1350
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1100
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 8999
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 8999
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 26997
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -36837
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3335.6667
$ws.Range("I70").Value = 3335.6667
$ws.Range("K70").Value = 3335.6667
$ws.Range("M70").Value = -3065.6667
$ws.Range("H73").Value = 3335.6667
$ws.Range("I73").Value = 3335.6667
$ws.Range("K73").Value = 3335.6667
$ws.Range("M73").Value = -2399.6667
$ws.Range("H80").Value = 2775.889
$ws.Range("I80").Value = 2436.8
$ws.Range("K80").Value = 2436.8
$ws.Range("M80").Value = -1438.8
$ws.Range("H83").Value = 2775.889
$ws.Range("I83").Value = 2436.8
$ws.Range("K83").Value = 12184
$ws.Range("M83").Value = -7192
$ws.Range("H102").Value = 2092.75
$ws.Range("J102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("N102").Value = -5244
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7818.5557
$ws.Range("I7").Value = 8248.5
$ws.Range("K7").Value = 8248.5
$ws.Range("M7").Value = -8136.5
$ws.Range("H40").Value = 4750.25
$ws.Range("I40").Value = 4200.4
$ws.Range("K40").Value = 4200.4
$ws.Range("M40").Value = -4064.4
$ws.Range("H46").Value = 1693.7778
$ws.Range("I46").Value = 1116.3334
$ws.Range("K46").Value = 1116.3334
$ws.Range("M46").Value = -928.3334
$ws.Range("H61").Value = 1571.2858
$ws.Range("H113").Value = 1571.2858
$ws.Range("H122").Value = 6108.423
$ws.Range("I122").Value = 4880.067
$ws.Range("K122").Value = 14640.201
$ws.Range("M122").Value = -12190.201
$ws.Range("H126").Value = 7818.5557
$ws.Range("I126").Value = 8248.5
$ws.Range("K126").Value = 24745.5
$ws.Range("M126").Value = -22275.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50629
$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -52183
$ws.Range("H136").Value = 2038.8667
$ws.Range("I136").Value = 2038.8667
$ws.Range("K136").Value = 6116.6001
$ws.Range("M136").Value = -3566.6001
